# Auto-generated cell update script based on the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Text)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

Set-TextCell $ws.Range("D2") "26.337.39"
$ws.Range("E2").Value = "  +0.96%  "
Set-TextCell $ws.Range("D3") "1.666.70"
$ws.Range("E3").Value = "  +0.87%  "
Set-TextCell $ws.Range("D4") "1.010"
$ws.Range("E4").Value = "  +0.88%  "
Set-TextCell $ws.Range("D5") "219.22"
$ws.Range("E5").Value = "  +0.81%  "
Set-TextCell $ws.Range("D6") "0.5345"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("E7").Value = "  +0.84%  "
Set-TextCell $ws.Range("D8") "0.2665"
$ws.Range("E8").Value = "  +2.58%  "
Set-TextCell $ws.Range("D9") "0.06402"
$ws.Range("E9").Value = "  +1.20%  "
Set-TextCell $ws.Range("D10") "20.93"
$ws.Range("E10").Value = "  +2.72%  "
Set-TextCell $ws.Range("D11") "0.07852"
$ws.Range("E11").Value = "  +0.83%  "
Set-TextCell $ws.Range("D12") "4.562"
$ws.Range("E12").Value = "  +1.10%  "
Set-TextCell $ws.Range("D13") "1.663.52"
$ws.Range("E13").Value = "  +0.81%  "
Set-TextCell $ws.Range("D14") "1.894.78"
$ws.Range("E14").Value = "  +0.84%  "
Set-TextCell $ws.Range("D15") "0.5542"
$ws.Range("E15").Value = "  +0.86%  "
Set-TextCell $ws.Range("D16") "0.0₅8200"
$ws.Range("E16").Value = "  -0.06%  "
Set-TextCell $ws.Range("D17") "65.89"
$ws.Range("E17").Value = "  +0.45%  "
Set-TextCell $ws.Range("D18") "26.358.18"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  +0.85%  "
Set-TextCell $ws.Range("D20") "4.689"
$ws.Range("E20").Value = "  +2.36%  "
Set-TextCell $ws.Range("D21") "193.91"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  +2.17%  "
Set-TextCell $ws.Range("D23") "6.047"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  +0.90%  "
Set-TextCell $ws.Range("D25") "146.35"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("E26").Value = "  -0.67%  "
Set-TextCell $ws.Range("D27") "7.214"
$ws.Range("E27").Value = "  -0.30%  "
Set-TextCell $ws.Range("D28") "16.12"
$ws.Range("E28").Value = "  +0.42%  "
Set-TextCell $ws.Range("D29") "1.499"
$ws.Range("E29").Value = "  +4.68%  "
Set-TextCell $ws.Range("D30") "0.05869"
$ws.Range("E30").Value = "  +0.45%  "
Set-TextCell $ws.Range("D31") "1.284"
$ws.Range("E31").Value = "  +0.83%  "
Set-TextCell $ws.Range("D32") "3.603"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("E33").Value = "  +0.51%  "
Set-TextCell $ws.Range("D34") "1.604"
$ws.Range("E34").Value = "  +1.25%  "
Set-TextCell $ws.Range("D35") "0.9718"
$ws.Range("E35").Value = "  +2.84%  "
Set-TextCell $ws.Range("D36") "2.827"
$ws.Range("E36").Value = "  +1.67%  "
Set-TextCell $ws.Range("D37") "2.421"
$ws.Range("E37").Value = "  +0.45%  "
Set-TextCell $ws.Range("D38") "0.5836"
$ws.Range("E38").Value = "  +1.71%  "
Set-TextCell $ws.Range("D39") "0.01608"
$ws.Range("E39").Value = "  -0.04%  "
Set-TextCell $ws.Range("D40") "0.8626"
$ws.Range("E40").Value = "  +2.30%  "
Set-TextCell $ws.Range("D41") "1.066.03"
Set-TextCell $ws.Range("D42") "5.839"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell $ws.Range("D43") "1.011"
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell $ws.Range("D44") "104.95"
$ws.Range("E44").Value = "  +0.54%  "
Set-TextCell $ws.Range("D45") "1.805.76"
$ws.Range("E45").Value = "  +0.63%  "
Set-TextCell $ws.Range("D46") "57.85"
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell $ws.Range("D47") "1.013"
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell $ws.Range("D48") "0.0₈104"
$ws.Range("E48").Value = "  -7.09%  "
Set-TextCell $ws.Range("D49") "0.4390"
$ws.Range("E49").Value = "  +1.37%  "
Set-TextCell $ws.Range("D50") "8.011"
$ws.Range("E50").Value = "  +2.73%  "
$ws.Range("E51").Value = "  +0.40%  "
